$wb = $excel.ActiveWorkbook

# --- Update header row on the "CypherOutput" sheet ---
$wsData = $wb.Worksheets.Item("CypherOutput")
$wsData.Range("A1").Value = "Case ID"
$wsData.Range("B1").Value = "Study Code"
$wsData.Range("C1").Value = "Study Type"
$wsData.Range("D1").Value = "Breed"
$wsData.Range("E1").Value = "Diagnosis"
$wsData.Range("F1").Value = "Stage of Disease"
$wsData.Range("G1").Value = "Age"
$wsData.Range("H1").Value = "Sex"
$wsData.Range("I1").Value = "Neutered Status"

# --- Update the Cypher query and output path on the "Message" sheet ---
$wsMsg = $wb.Worksheets.Item("Message")
$wsMsg.Range("A8").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN [''COTC007B''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`'
$wsMsg.Range("A10").Value = 'C:\Users\radhakrishnang2\Desktop\DataCommons_Automation\CTDC_Automation\OutputFiles\TC01_Canine_Filter_Study-COTB_Neo4jData.xlsx'
